$d = $word.ActiveDocument

$pairs = @(
    @("513÷8=", "216÷6="),
    @("810÷4=", "618÷5="),
    @("876÷4=", "110÷8="),
    @("307÷9=", "403÷2="),
    @("238÷2=", "362÷6="),
    @("531÷9=", "546÷9="),
    @("734÷6=", "309÷8="),
    @("184÷3=", "314÷5="),
    @("981÷3=", "954÷7="),
    @("898÷6=", "465÷2="),
    @("170÷5=", "369÷9="),
    @("711÷3=", "524÷6="),
    @("550÷6=", "371÷3="),
    @("143÷9=", "441÷4="),
    @("924÷8=", "424÷9="),
    @("327÷9=", "451÷6="),
    @("247÷2=", "259÷2="),
    @("840÷9=", "586÷6="),
    @("513÷6=", "782÷7="),
    @("805÷9=", "373÷9="),
    @("200÷7=", "117÷5="),
    @("155÷7=", "179÷2="),
    @("400÷2=", "404÷2="),
    @("191÷2=", "931÷5="),
    @("627÷8=", "164÷8=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
